# Apply edits: remap each row's target_cat/category/cond_cat/correct_answer/stimulus/
# conceptual/perceptual/typicality/n/p_*/r_* values from a different source row (per the
# commit's stimulus-set shuffle), and renumber trial_total (F) sequentially from 82.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns: Row, F(trial_total), G(target_cat), H(category), I(cond_cat), K(correct_answer),
#          L(stimulus), M(conceptual), N(perceptual), O(typicality), P(n),
#          Q(p_typicality), R(p_conceptual), S(p_perceptual), T(r_typicality),
#          U(r_conceptual), V(r_perceptual)
$rowsData = @(
    @(2, 82, "living_rooms", "living_rooms", "target", "j", "stimuli/img_njhlh.png", 59.74418604651163, 41.51162790697674, 50.62790697674419, 43, 4, 4, 4, 4, 4, 4),
    @(3, 83, "living_rooms", "living_rooms", "target", "j", "stimuli/img_j73b6.png", 21.5609756097561, 19.90243902439024, 20.73170731707317, 41, 1, 1, 1, 2, 1, 2),
    @(4, 84, "living_rooms", "living_rooms", "target", "j", "stimuli/img_9684y.png", 77.95744680851064, 56.70212765957447, 67.32978723404256, 47, 7, 7, 7, 6, 7, 7),
    @(5, 85, "living_rooms", "living_rooms", "target", "j", "stimuli/img_vg73h.png", 87.7, 72.4, 80.05000000000001, 50, 10, 10, 10, 10, 10, 10),
    @(6, 86, "living_rooms", "living_rooms", "target", "j", "stimuli/img_3ze38.png", 35.53191489361702, 28.4468085106383, 31.98936170212766, 47, 2, 2, 2, 3, 2, 3),
    @(7, 87, "living_rooms", "living_rooms", "target", "j", "stimuli/img_bg264.png", 87.9047619047619, 71.5, 79.70238095238095, 42, 10, 10, 10, 9, 10, 9),
    @(8, 88, "living_rooms", "kitchens", "distractor", "f", "stimuli/img_z10c7.png", 63.45945945945946, 37.21621621621622, 50.33783783783784, 37, 3, 3, 3, 3, 3, 3),
    @(9, 89, "living_rooms", "bedrooms", "distractor", "f", "stimuli/img_die1d.png", 75.42857142857143, 53.30952380952381, 64.36904761904762, 42, 6, 6, 6, 6, 6, 6),
    @(10, 90, "living_rooms", "kitchens", "distractor", "f", "stimuli/img_91lbx.png", 63.97058823529412, 38.8235294117647, 51.39705882352941, 34, 3, 3, 3, 3, 3, 3),
    @(11, 91, "living_rooms", "bedrooms", "distractor", "f", "stimuli/img_gxm46.png", 74.78378378378379, 54, 64.3918918918919, 37, 6, 6, 6, 6, 6, 6),
    @(12, 92, "living_rooms", "bedrooms", "distractor", "f", "stimuli/img_ds88w.png", 61.425, 38.2, 49.8125, 40, 3, 3, 3, 3, 3, 3),
    @(13, 93, "living_rooms", "living_rooms", "target", "j", "stimuli/img_tujn3.png", 81.4090909090909, 62.52272727272727, 71.9659090909091, 44, 8, 8, 8, 7, 8, 7),
    @(14, 94, "living_rooms", "living_rooms", "target", "j", "stimuli/img_tv8e2.png", 71.93023255813954, 50.25581395348837, 61.09302325581395, 43, 6, 6, 6, 6, 6, 6),
    @(15, 95, "living_rooms", "living_rooms", "target", "j", "stimuli/img_kq9s9.png", 62.30232558139535, 39.97674418604651, 51.13953488372093, 43, 4, 4, 4, 5, 5, 4),
    @(16, 96, "living_rooms", "bedrooms", "distractor", "f", "stimuli/img_qgbyn.png", 65.08108108108108, 40.10810810810811, 52.5945945945946, 37, 4, 4, 4, 4, 4, 4),
    @(17, 97, "living_rooms", "living_rooms", "target", "j", "stimuli/img_2qhro.png", 81.73809523809524, 62.73809523809524, 72.23809523809524, 42, 8, 8, 8, 8, 8, 8),
    @(18, 98, "living_rooms", "living_rooms", "target", "j", "stimuli/img_bf8nx.png", 86.63414634146342, 66.63414634146342, 76.63414634146342, 41, 9, 9, 9, 8, 9, 8),
    @(19, 99, "living_rooms", "kitchens", "distractor", "f", "stimuli/img_k0ze7.png", 58.82142857142857, 37.46428571428572, 48.14285714285714, 28, 3, 3, 3, 3, 3, 3),
    @(20, 100, "living_rooms", "bedrooms", "distractor", "f", "stimuli/img_jr212.png", 59.48571428571429, 39.08571428571429, 49.28571428571429, 35, 3, 3, 3, 3, 3, 3),
    @(21, 101, "living_rooms", "kitchens", "distractor", "f", "stimuli/img_9mcah.png", 59.59375, 36.75, 48.171875, 32, 3, 3, 3, 3, 3, 3),
    @(22, 102, "living_rooms", "bedrooms", "distractor", "f", "stimuli/img_2dnfy.png", 59.97297297297298, 37.56756756756756, 48.77027027027027, 37, 3, 3, 3, 3, 3, 3),
    @(23, 103, "living_rooms", "bedrooms", "distractor", "f", "stimuli/img_ccn2w.png", 65.7872340425532, 40.31914893617022, 53.05319148936171, 47, 4, 4, 4, 4, 4, 4),
    @(24, 104, "living_rooms", "kitchens", "distractor", "f", "stimuli/img_1r08n.png", 72.34285714285714, 49, 60.67142857142857, 35, 5, 5, 5, 5, 5, 5),
    @(25, 105, "living_rooms", "bedrooms", "distractor", "f", "stimuli/img_mdpr4.png", 74.04255319148936, 54.70212765957447, 64.37234042553192, 47, 6, 6, 6, 6, 6, 6),
    @(26, 106, "living_rooms", "living_rooms", "target", "j", "stimuli/img_6wne4.png", 25.16279069767442, 15, 20.08139534883721, 43, 1, 1, 1, 1, 2, 1),
    @(27, 107, "living_rooms", "living_rooms", "target", "j", "stimuli/img_99exx.png", 70.02272727272727, 51.88636363636363, 60.95454545454545, 44, 6, 6, 6, 5, 5, 6),
    @(28, 108, "living_rooms", "kitchens", "distractor", "f", "stimuli/img_2b8fp.png", 73.89189189189189, 51.45945945945946, 62.67567567567568, 37, 6, 6, 6, 6, 6, 6),
    @(29, 109, "living_rooms", "living_rooms", "target", "j", "stimuli/img_vnxft.png", 53.22727272727273, 34.84090909090909, 44.03409090909091, 44, 3, 3, 3, 4, 3, 4),
    @(30, 110, "living_rooms", "kitchens", "distractor", "f", "stimuli/img_3gm8h.png", 65.07894736842105, 43.92105263157895, 54.5, 38, 4, 4, 4, 4, 4, 4),
    @(31, 111, "living_rooms", "living_rooms", "target", "j", "stimuli/img_mdh76.png", 37.31914893617022, 25.12765957446809, 31.22340425531915, 47, 2, 2, 2, 2, 3, 2),
    @(32, 112, "living_rooms", "kitchens", "distractor", "f", "stimuli/img_x9ml4.png", 58.71428571428572, 38.46428571428572, 48.58928571428572, 28, 3, 3, 3, 3, 3, 3),
    @(33, 113, "living_rooms", "kitchens", "distractor", "f", "stimuli/img_di49w.png", 73.03125, 52.25, 62.640625, 32, 6, 6, 6, 6, 6, 6),
    @(34, 114, "living_rooms", "bedrooms", "distractor", "f", "stimuli/img_49h18.png", 59.28947368421053, 32.8421052631579, 46.06578947368421, 38, 3, 3, 3, 3, 3, 3),
    @(35, 115, "living_rooms", "living_rooms", "target", "j", "stimuli/img_qdln8.png", 85.51162790697674, 67.86046511627907, 76.68604651162791, 43, 9, 9, 9, 9, 9, 9),
    @(36, 116, "living_rooms", "living_rooms", "target", "j", "stimuli/img_sx68r.png", 54, 33.2051282051282, 43.6025641025641, 39, 3, 3, 3, 3, 4, 3),
    @(37, 117, "living_rooms", "living_rooms", "target", "j", "stimuli/img_di6f0.png", 94.04347826086956, 83.34782608695652, 88.69565217391303, 46, 10, 10, 10, 10, 10, 10),
    @(38, 118, "living_rooms", "living_rooms", "target", "j", "stimuli/img_xbtev.png", 13.68181818181818, 8.568181818181818, 11.125, 44, 1, 1, 1, 1, 1, 1),
    @(39, 119, "living_rooms", "kitchens", "distractor", "f", "stimuli/img_1r2ri.png", 62.44117647058823, 40.76470588235294, 51.60294117647059, 34, 3, 3, 3, 3, 3, 3),
    @(40, 120, "living_rooms", "living_rooms", "target", "j", "stimuli/img_r10cu.png", 78.52380952380952, 56.14285714285715, 67.33333333333333, 42, 7, 7, 7, 7, 7, 6),
    @(41, 121, "living_rooms", "bedrooms", "distractor", "f", "stimuli/img_gsfx4.png", 59.94736842105263, 33.81578947368421, 46.88157894736842, 38, 3, 3, 3, 3, 3, 3)
)

foreach ($r in $rowsData) {
    $rowNum = $r[0]
    $ws.Cells.Item($rowNum, 6).Value  = $r[1]   # F
    $ws.Cells.Item($rowNum, 7).Value  = $r[2]   # G
    $ws.Cells.Item($rowNum, 8).Value  = $r[3]   # H
    $ws.Cells.Item($rowNum, 9).Value  = $r[4]   # I
    $ws.Cells.Item($rowNum, 11).Value = $r[5]   # K
    $ws.Cells.Item($rowNum, 12).Value = $r[6]   # L
    $ws.Cells.Item($rowNum, 13).Value = $r[7]   # M
    $ws.Cells.Item($rowNum, 14).Value = $r[8]   # N
    $ws.Cells.Item($rowNum, 15).Value = $r[9]   # O
    $ws.Cells.Item($rowNum, 16).Value = $r[10]  # P
    $ws.Cells.Item($rowNum, 17).Value = $r[11]  # Q
    $ws.Cells.Item($rowNum, 18).Value = $r[12]  # R
    $ws.Cells.Item($rowNum, 19).Value = $r[13]  # S
    $ws.Cells.Item($rowNum, 20).Value = $r[14]  # T
    $ws.Cells.Item($rowNum, 21).Value = $r[15]  # U
    $ws.Cells.Item($rowNum, 22).Value = $r[16]  # V
}
